$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("Config")
$wsOptions = $wb.Worksheets.Item("ConfigOptions")

# --- ConfigOptions sheet: update the scroll/selection view state only ---
$wsOptions.Activate()
$wsOptions.Range("B25:B26").Select()

# --- Config sheet: update config values ---
$wsConfig.Activate()

# MasterFolder value updated to the new local OneDrive desktop path
$wsConfig.Range("B6").Value = "C:\Users\RollLe01\OneDrive - Reed Elsevier Group ICO Reed Elsevier Inc\Desktop\FLOBOT\"

# ZippedDirectory value updated to the new local OneDrive desktop path
$wsConfig.Range("B8").Value = "C:\Users\RollLe01\OneDrive - Reed Elsevier Group ICO Reed Elsevier Inc\Desktop\FLOBOT"

# Add "selector" mailto hyperlinks for the RecipientTo / RecipientCC-adjacent
# rows (B11 first, then B10, matching the order the links were inserted)
$wsConfig.Hyperlinks.Add($wsConfig.Range("B11"), "mailto:lester.rollan@lexisnexisrisk.com")
$wsConfig.Hyperlinks.Add($wsConfig.Range("B10"), "mailto:lester.rollan@lexisnexisrisk.com")

# Re-apply the existing Hyperlink cell style (copied from the already-styled
# B3 hyperlink cell) so B10/B11 reuse the workbook's existing Hyperlink xf
# instead of a freshly minted one.
$wsConfig.Range("B3").Copy()
$wsConfig.Range("B10:B11").PasteSpecial(-4122)

# Update the active selection to reflect the cell of interest after the edit
$wsConfig.Range("B10").Select()

Write-Host "edit complete"
